$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("g10.1")

$ws.Cells.Item(2, 2).Value = -6.943213131802861
$ws.Cells.Item(2, 3).Value = 7.980742702876276
$ws.Cells.Item(2, 4).Value = 3.559178325481649
$ws.Cells.Item(3, 2).Value = -3.033399858582009
$ws.Cells.Item(3, 3).Value = 4.325288184342346
$ws.Cells.Item(3, 4).Value = -5.228005074708375
$ws.Cells.Item(4, 2).Value = -0.2288600019107689
$ws.Cells.Item(4, 3).Value = 4.090352838664701
$ws.Cells.Item(4, 4).Value = 1.61985019946338
$ws.Cells.Item(5, 2).Value = 1.24772536898583
$ws.Cells.Item(5, 3).Value = -1.19499485174599
$ws.Cells.Item(5, 4).Value = 8.071533547379129
$ws.Cells.Item(6, 2).Value = -5.140084324314231
$ws.Cells.Item(6, 3).Value = -3.652084222111918
$ws.Cells.Item(6, 4).Value = 0.1774081972812258
$ws.Cells.Item(7, 2).Value = -3.397932324802488
$ws.Cells.Item(7, 3).Value = 0.583525770808202
$ws.Cells.Item(7, 4).Value = 0.7598352624477389
$ws.Cells.Item(8, 2).Value = -3.18596839631059
$ws.Cells.Item(8, 3).Value = -1.091482996358195
$ws.Cells.Item(8, 4).Value = -1.922784329967397
$ws.Cells.Item(9, 2).Value = 2.972239650855424
$ws.Cells.Item(9, 3).Value = 1.038283775507809
$ws.Cells.Item(9, 4).Value = 9.604857944187906
$ws.Cells.Item(10, 2).Value = -13.97255252459051
$ws.Cells.Item(10, 3).Value = -3.906522121507139
$ws.Cells.Item(10, 4).Value = -12.17208917510788
$ws.Cells.Item(11, 2).Value = -11.05852305620129
$ws.Cells.Item(11, 3).Value = 15.38588480891123
$ws.Cells.Item(11, 4).Value = -14.12283002730265
$ws.Cells.Item(12, 2).Value = -4.755431186326897
$ws.Cells.Item(12, 3).Value = 14.56074928380076
$ws.Cells.Item(12, 4).Value = -12.94489690617815
$ws.Cells.Item(13, 2).Value = -6.003665244428714
$ws.Cells.Item(13, 3).Value = 7.569339803891406
$ws.Cells.Item(13, 4).Value = -7.531773508934014
